$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 32 / 33: the match details (columns F:V) were swapped between the
#    two rows (the fixed per-row columns A:E — index/country/tournament/
#    season/kickoff — stay put).
# ---------------------------------------------------------------------------
$row32 = $ws.Range("F32:V32").Value()
$row33 = $ws.Range("F33:V33").Value()
$ws.Range("F32:V32").Value = $row33
$ws.Range("F33:V33").Value = $row32

# ---------------------------------------------------------------------------
# 2) Rows 37 / 38: same kind of swap.
# ---------------------------------------------------------------------------
$row37 = $ws.Range("F37:V37").Value()
$row38 = $ws.Range("F38:V38").Value()
$ws.Range("F37:V37").Value = $row38
$ws.Range("F38:V38").Value = $row37

# ---------------------------------------------------------------------------
# 3) Append four new match rows (68-71) at the bottom of the sheet.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=68; A=67; E=45190.91666666666; F="Belgrano";          G=3; H="Platense";         I=0;
       J=2.17; K="17/09/2023 21:11"; L=2.28; M="21/09/2023 21:52";
       N=3.09; O="17/09/2023 21:11"; P=2.76; Q="21/09/2023 21:52";
       R=4.01; S="17/09/2023 21:11"; T=4.28; U="21/09/2023 21:52";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/ca-belgrano-de-cordoba-platense/pneGCCfT/" },

    @{ Row=69; A=68; E=45190.97916666666; F="Union de Santa Fe"; G=0; H="Godoy Cruz";        I=0;
       J=2.08; K="17/09/2023 21:11"; L=2.31; M="21/09/2023 23:26";
       N=3.27; O="17/09/2023 21:11"; P=3;    Q="21/09/2023 23:24";
       R=4.05; S="17/09/2023 21:11"; T=3.75; U="21/09/2023 23:26";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/union-de-santa-fe-godoy-cruz/Quwq8W2p/" },

    @{ Row=70; A=69; E=45191.01041666666; F="Argentinos Jrs";    G=3; H="Talleres Cordoba";  I=1;
       J=2.21; K="17/09/2023 01:12"; L=1.91; M="22/09/2023 00:14";
       N=3.11; O="17/09/2023 01:12"; P=3.51; Q="22/09/2023 00:14";
       R=3.71; S="17/09/2023 01:12"; T=4.43; U="22/09/2023 00:14";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/argentinos-jrs-talleres-cordoba/WI1mJdBo/" },

    @{ Row=71; A=70; E=45191.08333333334; F="River Plate";      G=1; H="Atl. Tucuman";      I=0;
       J=1.39; K="16/09/2023 23:12"; L=1.45; M="22/09/2023 01:55";
       N=4.72; O="16/09/2023 23:12"; P=4.62; Q="22/09/2023 01:58";
       R=7.33; S="16/09/2023 23:12"; T=7.56; U="22/09/2023 01:58";
       V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/river-plate-atl-tucuman/IyeyCu4I/" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Clone the formatting (not the values) of the previous row so the new
    # row picks up the same per-column styling (bold/bordered index column,
    # datetime-formatted kickoff column, plain cells elsewhere).
    $ws.Range("A" + ($rowNum - 1) + ":V" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum + ":V" + $rowNum).PasteSpecial(-4122)

    # B:D (pais/torneio/temporada) are identical text on every row of this
    # sheet - copy them from row 2 so "2023" etc. stay text, not numbers.
    $ws.Range("B2:D2").Copy()
    $ws.Range("B" + $rowNum + ":D" + $rowNum).PasteSpecial()

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("E" + $rowNum).Value = $r.E
    $ws.Range("F" + $rowNum).Value = $r.F
    $ws.Range("G" + $rowNum).Value = $r.G
    $ws.Range("H" + $rowNum).Value = $r.H
    $ws.Range("I" + $rowNum).Value = $r.I
    $ws.Range("J" + $rowNum).Value = $r.J
    $ws.Range("K" + $rowNum).Value = $r.K
    $ws.Range("L" + $rowNum).Value = $r.L
    $ws.Range("M" + $rowNum).Value = $r.M
    $ws.Range("N" + $rowNum).Value = $r.N
    $ws.Range("O" + $rowNum).Value = $r.O
    $ws.Range("P" + $rowNum).Value = $r.P
    $ws.Range("Q" + $rowNum).Value = $r.Q
    $ws.Range("R" + $rowNum).Value = $r.R
    $ws.Range("S" + $rowNum).Value = $r.S
    $ws.Range("T" + $rowNum).Value = $r.T
    $ws.Range("U" + $rowNum).Value = $r.U
    $ws.Range("V" + $rowNum).Value = $r.V
}
